$d = $word.ActiveDocument

# Change 1: merge three runs (remove the run-split around "data de nascimento")
# by simply replacing the old concatenated text pattern with itself collapsed —
# achieved via Find/Replace across the run boundary.
$d.Content.Find.Execute(
    "idade, data de nascimento e o nome.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "idade, data de nascimento e o nome.", 2) | Out-Null

# Change 2: extend the investment-simulation sentence with new text.
$d.Content.Find.Execute(
    "tesouro etc).",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "tesouro etc) que deve armazenar a descrição, rendimento e o período que vai render.", 2) | Out-Null
